# Adds two new worksheets ("11" and "12") to the workbook. Each sheet holds
# a header-less grid of pendulum-observation rows: column A carries the
# textual "conditional(...)"/"add(...)" expression tree identifying that
# row's formula, and columns B:Q carry the 16 numeric trajectory samples
# produced by evaluating it. The two sheets are appended after the existing
# ones (Sheet1, "9.81"); the last one added ("12") ends up selected, which
# is what drives the workbook's activeTab changing from 1 to 3 (0-based:
# Sheet1=0, 9.81=1, 11=2, 12=3).

$wb = $excel.ActiveWorkbook

$data11 = @(
    @('conditional(add(y, y), add(conditional(x, conditional(conditional(y, x), add(vel, vel))), add(add(x, add(vel, add(conditional(conditional(conditional(x, vel), conditional(vel, vel)), conditional(x, x)), add(x, x)))), add(y, y))))', -1244.91, -1182.54, -935.77, -1437.9, -1662.46, -1600.89, -1502.35, -1631.74, -1815.57, -1902.44, -1968.06, -2022.6, -2044.17, -2143.11, -2195.58, -2211.44),
    @('conditional(add(y, y), add(add(add(add(x, y), x), add(add(x, add(x, x)), y)), vel))', -1107.55, -1216.59, -1267.9, -1468.39, -1689.83, -1664.02, -1594.55, -1514.33, -1803.23, -1920.47, -1941.31, -2002.31, -2066.33, -2166.54, -2174.14, -2212.99),
    @('add(conditional(add(y, y), add(add(add(conditional(conditional(x, y), conditional(vel, conditional(x, add(x, vel)))), add(add(x, add(vel, x)), vel)), add(x, y)), add(add(x, x), x))), conditional(conditional(y, y), add(add(add(add(conditional(x, conditional(conditional(y, conditional(add(x, y), add(y, y))), vel)), add(x, vel)), add(x, conditional(conditional(x, add(y, y)), conditional(conditional(conditional(vel, x), vel), add(y, x))))), add(add(x, x), conditional(add(conditional(x, x), y), add(x, conditional(y, y))))), y)))', -1264.06, -1177.72, -1114.32, -1204.95, -1357.91, -1446.38, -1317.76, -1331.27, -1772.33, -1930.66, -1929, -1998.95, -2029.1, -2118.43, -2175.85, -2207.69),
    @('conditional(add(y, y), add(conditional(conditional(conditional(conditional(x, vel), conditional(add(vel, x), conditional(y, x))), conditional(conditional(add(conditional(x, x), vel), x), vel)), y), add(add(x, y), add(add(vel, add(x, add(y, x))), conditional(x, x)))))', -1210.09, -1105.84, -1260.36, -1531.08, -1954.74, -2056.67, -1899.43, -1895.74, -1914.2, -1939.47, -2033.78, -2032.25, -2092.22, -2143.18, -2179.92, -2239.04),
    @('add(conditional(y, add(x, vel)), conditional(add(y, y), add(add(add(add(x, add(conditional(add(y, x), add(y, x)), x)), x), add(vel, conditional(y, y))), conditional(add(conditional(x, conditional(vel, add(y, conditional(add(vel, y), x)))), x), conditional(vel, add(add(vel, x), conditional(add(y, y), conditional(add(conditional(y, vel), conditional(x, vel)), add(x, vel)))))))))', -1129.42, -1444.57, -1495.71, -1581.78, -1615.55, -1836.78, -1768.99, -1773.55, -1525.25, -1854.48, -1951.71, -2010.4, -2041.71, -2101.52, -2180.93, -2216.04),
    @('conditional(add(y, add(y, y)), add(add(add(y, add(conditional(conditional(x, x), add(add(add(add(add(add(x, y), add(conditional(vel, y), add(x, add(add(vel, add(vel, x)), conditional(y, y))))), conditional(y, x)), add(x, x)), add(add(y, conditional(x, conditional(add(y, vel), conditional(y, y)))), conditional(x, x))), conditional(conditional(add(vel, x), conditional(x, x)), add(y, y)))), add(x, y))), add(x, x)), add(x, add(add(vel, x), conditional(conditional(y, x), y)))))', -1207.71, -1188.63, -1215.32, -1169.53, -1414.98, -1520.92, -1502.04, -1596.81, -1809.02, -1919.52, -1978.17, -2033.31, -2081.88, -2159.23, -2167.15, -2233.17),
    @('conditional(add(conditional(y, y), add(y, y)), add(add(add(add(add(x, x), add(vel, x)), add(add(conditional(x, vel), x), add(x, add(add(add(add(add(vel, x), add(x, y)), add(add(x, vel), add(vel, conditional(vel, vel)))), add(add(add(add(add(add(add(x, x), add(y, y)), add(x, add(vel, add(add(y, vel), add(x, x))))), add(add(x, x), add(x, x))), x), x), x)), add(x, x))))), add(add(x, x), x)), x))', -1216.19, -1202.22, -1449.54, -1371.71, -1600.73, -1748.33, -1347.65, -1623.82, -1790.96, -1886.58, -1928.33, -2006.09, -2095.52, -2128.8, -2167.31, -2185.62),
    @('conditional(add(conditional(x, add(x, conditional(conditional(x, x), add(conditional(conditional(y, add(y, vel)), conditional(add(add(x, add(conditional(conditional(x, conditional(x, add(vel, x))), add(add(conditional(add(vel, x), conditional(vel, y)), x), add(y, vel))), y)), y), y)), add(x, add(y, y)))))), add(y, y)), add(add(add(add(x, y), y), x), add(vel, x)))', -1190.17, -1110.45, -1096.64, -1139.21, -1307.03, -1337.22, -1221.97, -1033.82, -1638.99, -1764.54, -1829.76, -1883.99, -1982.36, -2024.19, -2086.97, -2167.99),
    @('conditional(add(y, y), add(add(x, y), add(add(conditional(conditional(conditional(x, add(y, add(vel, x))), add(add(x, x), x)), add(x, add(y, y))), add(conditional(y, y), vel)), add(conditional(x, x), conditional(x, x)))))', -1371.04, -1240.46, -1039.53, -1269.1, -1312.16, -1241.11, -1528.74, -1529.18, -1731.15, -1914.64, -1938.13, -1977.21, -2048.92, -2104.07, -2181.33, -2198.5),
    @('conditional(add(add(y, y), y), add(add(add(x, vel), add(add(add(add(add(x, add(vel, add(x, add(x, x)))), add(add(y, add(y, vel)), add(x, add(add(conditional(x, x), conditional(y, y)), conditional(add(vel, x), conditional(vel, conditional(add(y, vel), add(y, vel)))))))), add(add(x, add(x, conditional(x, vel))), conditional(x, conditional(x, add(x, add(x, x)))))), add(x, x)), add(x, add(add(conditional(x, x), conditional(y, vel)), conditional(conditional(vel, vel), vel))))), add(add(x, add(x, x)), conditional(x, add(y, y)))))', -1305.89, -1331.15, -1380.47, -1367.48, -1657.79, -1953.13, -2075.77, -2156.55, -2034.84, -1982.33, -2060.46, -2128.71, -2079.07, -2135, -2221.96, -2239.86),
)

$data12 = @(
    @('conditional(add(y, y), add(add(add(y, x), conditional(conditional(y, y), y)), add(vel, add(add(add(add(x, x), conditional(conditional(x, x), conditional(x, vel))), conditional(x, add(add(conditional(conditional(vel, y), add(x, conditional(x, x))), add(add(add(conditional(add(x, vel), add(conditional(y, y), vel)), x), add(x, y)), add(vel, add(x, x)))), add(add(x, conditional(add(y, add(add(y, add(x, x)), y)), add(vel, vel))), conditional(x, x))))), x))))', -1153.42, -1114.02, -1031.71, -1182.4, -1292.69, -1227.98, -1482.45, -1662.67, -1830.02, -1909.82, -1901.46, -1980.47, -2003.76, -2092.23, -2175.86, -2190.39),
    @('conditional(add(conditional(x, add(conditional(add(y, conditional(conditional(vel, vel), add(vel, x))), add(add(add(x, vel), add(x, add(x, vel))), y)), add(add(y, x), conditional(x, vel)))), add(y, y)), add(add(x, vel), add(y, add(x, y))))', -1785.29, -1925.46, -2300.69, -2544.62, -2544.66, -2510.69, -2437.01, -2344.5, -2246.4, -2141.51, -2028.35, -1587.12, -1771.25, -1807.53, -1769.38, -2025.53),
    @('conditional(add(add(y, y), conditional(y, y)), add(add(add(conditional(x, vel), y), add(add(vel, x), add(conditional(y, y), conditional(x, x)))), add(add(y, x), add(conditional(x, x), conditional(x, x)))))', -1265.81, -1212.77, -1329.2, -1269.77, -1359.98, -1512.32, -1464.49, -1516.71, -1775.76, -1927.53, -1969.96, -2012.38, -2058.7, -2094.39, -2188.33, -2194.04),
    @('add(conditional(conditional(x, vel), add(x, x)), conditional(add(y, y), add(add(vel, add(conditional(conditional(x, y), conditional(conditional(conditional(x, add(conditional(vel, add(add(vel, vel), conditional(x, y))), add(vel, conditional(x, conditional(conditional(x, vel), add(x, x)))))), add(vel, x)), conditional(y, add(conditional(y, x), conditional(add(vel, conditional(x, y)), x))))), add(add(x, y), add(x, add(x, x))))), y)))', -2304.58, -2348.73, -2331.56, -2330.86, -2368.74, -2266.57, -2233.8, -1954.23, -1878.25, -2144.08, -2053.36, -1662.16, -1755.38, -1833.93, -1849.48, -2033.63),
    @('add(conditional(x, add(vel, add(add(y, x), conditional(conditional(vel, x), add(x, y))))), conditional(add(y, add(y, y)), add(add(add(add(add(x, conditional(conditional(x, add(conditional(x, x), conditional(conditional(x, vel), add(conditional(x, y), add(x, x))))), add(conditional(y, vel), vel))), add(add(x, x), add(x, y))), add(x, conditional(add(x, y), conditional(x, x)))), add(x, add(add(conditional(x, x), add(vel, add(x, y))), vel))), x)))', -1086.21, -1016.12, -821.5700000000001, -968.16, -1222.27, -1271.72, -1160.51, -1040.63, -596.37, -1131.75, -1381.04, -1574.01, -1741.91, -1808.34, -1992.16, -2056.4),
    @('add(conditional(x, conditional(vel, conditional(conditional(conditional(add(x, vel), conditional(conditional(add(vel, x), add(x, conditional(conditional(x, y), conditional(y, y)))), vel)), conditional(add(x, conditional(add(x, vel), add(vel, vel))), y)), add(conditional(y, y), add(y, add(vel, conditional(y, vel))))))), conditional(add(y, y), add(add(y, add(add(vel, x), add(x, x))), add(x, y))))', -2124.17, -2086.43, -2219.44, -2489.29, -2528.56, -2514.32, -2445.12, -2356.86, -2171.2, -2067.86, -1834.06, -1925.51, -1949.95, -1948.11, -1963.7, -2116.6),
    @('conditional(add(y, add(y, y)), add(add(add(x, add(add(add(conditional(conditional(x, add(conditional(x, x), conditional(conditional(add(vel, y), x), y))), conditional(x, y)), add(add(y, x), conditional(conditional(conditional(vel, conditional(y, conditional(y, y))), conditional(add(x, y), conditional(x, x))), conditional(add(vel, y), conditional(add(x, y), conditional(conditional(conditional(conditional(add(x, vel), add(conditional(x, x), conditional(x, x))), vel), conditional(vel, y)), vel)))))), x), vel)), add(add(x, y), add(x, y))), vel))', -1626.2, -1429.15, -1244.09, -1171.91, -1222.53, -1337.19, -1254.79, -1049.36, -1665.77, -1762.01, -1911.21, -1935.97, -2021.68, -2079.35, -2116.13, -2145.25),
    @('add(conditional(x, conditional(vel, add(conditional(conditional(conditional(y, x), conditional(conditional(conditional(conditional(y, x), conditional(conditional(conditional(add(x, y), vel), add(x, y)), add(x, add(conditional(conditional(vel, x), add(x, y)), conditional(y, vel))))), conditional(add(add(conditional(conditional(y, y), conditional(y, vel)), vel), y), y)), add(x, add(conditional(conditional(vel, x), add(x, y)), conditional(y, vel))))), conditional(add(add(conditional(conditional(y, y), conditional(y, vel)), vel), y), y)), add(add(add(x, add(x, conditional(conditional(vel, vel), add(add(y, add(add(y, x), conditional(x, vel))), vel)))), conditional(conditional(conditional(add(x, x), conditional(y, x)), y), conditional(conditional(x, conditional(y, conditional(add(vel, vel), add(conditional(add(y, x), add(y, y)), conditional(conditional(x, y), x))))), add(vel, add(add(conditional(vel, x), conditional(conditional(add(x, add(vel, x)), conditional(x, vel)), vel)), conditional(add(y, add(add(conditional(vel, conditional(conditional(x, y), conditional(y, vel))), add(vel, y)), add(y, add(conditional(add(y, y), y), conditional(y, vel))))), x)))))), x)))), conditional(add(y, y), add(add(add(add(x, add(x, y)), y), add(vel, conditional(x, x))), x)))', -2156.57, -2145.19, -2101.42, -2404.31, -2371.79, -2333.31, -2239.19, -2136.95, -2105.3, -2030.5, -1901.77, -1699.11, -1767.64, -1868.14, -1848.14, -2045.28),
)

$sheetsToAdd = @(
    @{ Name = "11"; Rows = $data11 },
    @{ Name = "12"; Rows = $data12 }
)

foreach ($entry in $sheetsToAdd) {
    $afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $afterSheet)
    $ws.Name = $entry.Name

    $rows = $entry.Rows
    for ($r = 0; $r -lt $rows.Count; $r++) {
        $row = $rows[$r]
        for ($c = 0; $c -lt $row.Count; $c++) {
            $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
        }
    }

    $ws.Range("A1").Select() | Out-Null
}

# Leave the last-added sheet ("12") as the active tab.
$ws.Activate()
